$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct typo in the recurring_challenge_text sentence template:
# "has been reported as challenge" -> "has been reported as a challenge"
$ws.Range("D4").Value = "**{challenge}** has been reported as a challenge for the **{goal}** team in each of the last **{challenge count} quarters**."

# Reflect the cursor/selection ending on the cell that was edited.
[void]$ws.Range("D4").Select()
